# Daily attendance processing - 2026-01-20 08:06:48
# Reorders the "Recorded By" (column G) author list so that the "System"
# entry swaps places with whichever entry it is paired with (the value
# immediately next to it among the first two comma-separated entries).
# Single-entry cells, and cells where "System" is not among the first two
# entries, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $parts = $value.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -ge 2 -and ($parts[0] -eq "System" -or $parts[1] -eq "System")) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp

        $newValue = [string]::Join(", ", $parts)
        $cell.Value2 = $newValue
    }
}
